# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 161
$wsExhibit.Range("F7").Value = 1600
$wsExhibit.Range("F10").Value = 1344
$wsExhibit.Range("F12").Value = 15
$wsExhibit.Range("F18").Value = 245
$wsExhibit.Range("F20").Value = 203

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 161
$wsAll.Range("F7").Value = 1600
$wsAll.Range("F11").Value = 1344
$wsAll.Range("F19").Value = 245
$wsAll.Range("F21").Value = 203
